# Updating Person Classification evaluation directory with latest confusion
# matrices and results: split the single-sheet results.xlsx into three
# per-model sheets (ArcFace, VGGFace, FaceNet512), each carrying
# Accuracy/Precision/Recall/F1-Score broken out by weighted/micro/macro
# averaging.

$wb = $excel.ActiveWorkbook

# The workbook starts with exactly one sheet ("Sheet1"); grab it before any
# structural edits happen.
$orig = $wb.Worksheets.Item(1)

# Duplicate the original sheet twice (placed right after it each time) so
# the two new per-model sheets inherit the same formatting/markup as the
# original, then rename everything into place.
$orig.Copy([System.Reflection.Missing]::Value, $orig)
$orig.Copy([System.Reflection.Missing]::Value, $orig)

$wb.Worksheets.Item("Sheet1 (2)").Name = "ArcFace"
$wb.Worksheets.Item("Sheet1 (3)").Name = "VGGFace"
$wb.Worksheets.Item("Sheet1").Name = "FaceNet512"

# Re-order into the final left-to-right tab order: ArcFace, VGGFace,
# FaceNet512.
$wb.Worksheets.Item("ArcFace").Move($wb.Worksheets.Item(1))
$wb.Worksheets.Item("VGGFace").Move([System.Reflection.Missing]::Value, $wb.Worksheets.Item("ArcFace"))

$arc = $wb.Worksheets.Item("ArcFace")
$vgg = $wb.Worksheets.Item("VGGFace")
$facenet = $wb.Worksheets.Item("FaceNet512")

$sheetsData = @(
    @{ ws = $arc;     accuracy = 0.903949;             precisionWeighted = 0.92196500000000003; precisionMacro = 0.62243300000000001; recallWeighted = 0.903949;             recallMacro = 0.58946399999999999; f1Weighted = 0.90428699999999995; f1Macro = 0.59957000000000005; selCell = "E6" },
    @{ ws = $vgg;     accuracy = 0.93489900000000004;  precisionWeighted = 0.94334499999999999; precisionMacro = 0.63424899999999995; recallWeighted = 0.93489900000000004;  recallMacro = 0.61476500000000001; f1Weighted = 0.93548500000000001; f1Macro = 0.62185999999999997; selCell = "E3" },
    @{ ws = $facenet; accuracy = 0.96904999999999997;  precisionWeighted = 0.975989;             precisionMacro = 0.65250600000000003; recallWeighted = 0.96904999999999997;  recallMacro = 0.64265799999999995; f1Weighted = 0.97201199999999999; f1Macro = 0.647204;                selCell = "E3" }
)

foreach ($d in $sheetsData) {
    $ws = $d.ws
    $acc = $d.accuracy

    $ws.Cells.Item(1,1).Value = "Metric"
    $ws.Cells.Item(1,2).Value = "Value (Weighted)"
    $ws.Cells.Item(1,3).Value = "Value (Micro)"
    $ws.Cells.Item(1,4).Value = "Value(Macro)"

    $ws.Cells.Item(2,1).Value = "Accuracy"
    $ws.Cells.Item(2,2).Value = $acc
    $ws.Cells.Item(2,3).Value = $acc
    $ws.Cells.Item(2,4).Value = $acc

    $ws.Cells.Item(3,1).Value = "Precision"
    $ws.Cells.Item(3,2).Value = $d.precisionWeighted
    $ws.Cells.Item(3,3).Value = $acc
    $ws.Cells.Item(3,4).Value = $d.precisionMacro

    $ws.Cells.Item(4,1).Value = "Recall"
    $ws.Cells.Item(4,2).Value = $d.recallWeighted
    $ws.Cells.Item(4,3).Value = $acc
    $ws.Cells.Item(4,4).Value = $d.recallMacro

    $ws.Cells.Item(5,1).Value = "F1-Score"
    $ws.Cells.Item(5,2).Value = $d.f1Weighted
    $ws.Cells.Item(5,3).Value = $acc
    $ws.Cells.Item(5,4).Value = $d.f1Macro

    $ws.Columns.Item(1).AutoFit() | Out-Null
    $ws.Columns.Item(2).AutoFit() | Out-Null
    $ws.Columns.Item(3).AutoFit() | Out-Null
    $ws.Columns.Item(4).AutoFit() | Out-Null

    $ws.Range($d.selCell).Select() | Out-Null
}

# FaceNet512 (the renamed original sheet) ends up the active/visible tab.
$facenet.Activate()
